$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2759
$ws.Range("L3").Value = 2802
$ws.Range("L4").Value = 749
$ws.Range("L5").Value = 161
$ws.Range("L6").Value = 2494
$ws.Range("L7").Value = 8965

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 35
$ws.Range("L7").Value = 301
$ws.Range("L8").Value = 566
$ws.Range("L11").Value = 158
$ws.Range("L15").Value = 65
$ws.Range("L18").Value = 63
$ws.Range("L19").Value = 249
$ws.Range("L20").Value = 230
$ws.Range("L24").Value = 20
$ws.Range("L29").Value = 481
$ws.Range("L31").Value = 86
$ws.Range("K31").Value = 327
$ws.Range("L33").Value = 409
$ws.Range("L34").Value = 54
$ws.Range("L36").Value = 123
$ws.Range("L37").Value = 329
$ws.Range("L41").Value = 41
$ws.Range("L42").Value = 298
$ws.Range("L43").Value = 69
$ws.Range("L44").Value = 67
$ws.Range("L47").Value = 68
$ws.Range("L50").Value = 48
$ws.Range("L52").Value = 180
$ws.Range("L54").Value = 183
$ws.Range("L60").Value = 55
$ws.Range("L61").Value = 11
$ws.Range("K63").Value = 158
$ws.Range("L63").Value = 30
$ws.Range("L64").Value = 55
$ws.Range("L67").Value = 331
$ws.Range("L69").Value = 24
$ws.Range("L73").Value = 76
$ws.Range("L74").Value = 9
$ws.Range("L76").Value = 112
$ws.Range("L79").Value = 240
$ws.Range("L83").Value = 211
$ws.Range("L88").Value = 113
$ws.Range("L89").Value = 114
$ws.Range("L94").Value = 109
$ws.Range("L96").Value = 88
$ws.Range("L97").Value = 80
$ws.Range("L99").Value = 149
$ws.Range("L101").Value = 8965

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 92
$ws.Range("L3").Value = 95
$ws.Range("L4").Value = 25
$ws.Range("L6").Value = 82
$ws.Range("L7").Value = 301

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 57
$ws.Range("L3").Value = 50
$ws.Range("L6").Value = 38
$ws.Range("L7").Value = 158

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 50
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 161
$ws.Range("L3").Value = 190
$ws.Range("L4").Value = 40
$ws.Range("L7").Value = 566

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L5").Value = 5
$ws.Range("L7").Value = 211

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 113
$ws.Range("L3").Value = 127
$ws.Range("L6").Value = 140
$ws.Range("L7").Value = 409

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 97
$ws.Range("L3").Value = 95
$ws.Range("L7").Value = 329

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 149

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L3").Value = 22
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 327
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 97
$ws.Range("L3").Value = 120
$ws.Range("L7").Value = 331

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 40
$ws.Range("L3").Value = 36
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 183

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 154
$ws.Range("L6").Value = 123
$ws.Range("L7").Value = 481

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 84
$ws.Range("L3").Value = 78
$ws.Range("L7").Value = 249

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 29
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 22
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 85
$ws.Range("L3").Value = 90
$ws.Range("L6").Value = 88
$ws.Range("L7").Value = 298

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("L3").Value = 5
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 87
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 240

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 72
$ws.Range("L3").Value = 70
$ws.Range("L7").Value = 230

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 123

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 30
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L3").Value = 17
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 80

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L6").Value = 35
$ws.Range("L7").Value = 113

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("Mount Greenwood")
$ws.Range("L6").Value = 4
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 9
